$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of the
# existing header row (bold font, thin border, centered/top aligned) by
# copying the format from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I (I0) and J (IF) for rows 2-17.
$data = @{
    2  = @(7, 9)
    3  = @(7, 7)
    4  = @(5, 5)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(8, 9)
    10 = @(3, 3)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(9, 9)
    15 = @(4, 5)
    16 = @(6, 6)
    17 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
